# The scraper's xpath expressions were adjusted, which dropped the
# "reviews_count" column (old column E) from the exported data. All of
# the columns that used to sit to its right (reviews_average, latitude,
# longitude, is_permanently_closed, gmaps_link, latest_review_date) shift
# one column to the left as a result.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E:E").Delete()
